# Apply the crypto-tracker refresh captured in the Nov 4 2024 GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "68.556.69"
$ws.Range("E2").Value = "  -0.01%  "

# Row 3
$ws.Range("D3").Value = "2.454.41"
$ws.Range("E3").Value = "  -0.24%  "

# Row 4
$ws.Range("E4").Value = "  -0.04%  "

# Row 5
$ws.Range("D5").Value = "'557.05"
$ws.Range("D5").Style = "Normal"  # keep as text, matching original (no numeric coercion)
$ws.Range("E5").Value = "  -1.11%  "

# Row 6
$ws.Range("D6").Value = "'160.37"
$ws.Range("D6").Style = "Normal"  # keep as text, matching original (no numeric coercion)
$ws.Range("E6").Value = "  -2.16%  "

# Row 7
$ws.Range("E7").Value = "  -0.03%  "

# Row 8
$ws.Range("E8").Value = "  -0.01%  "

# Row 9
$ws.Range("E9").Value = "  -1.51%  "

# Row 11
$ws.Range("B11").Value = "Cardano"
$ws.Range("C11").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D11").Value = "'0.330"
$ws.Range("D11").Style = "Normal"  # keep as text, matching original (no numeric coercion)
$ws.Range("E11").Value = "  -3.31%  "

# Row 12
$ws.Range("B12").Value = "Toncoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D12").Value = "'4.83"
$ws.Range("D12").Style = "Normal"  # keep as text, matching original (no numeric coercion)
$ws.Range("E12").Value = "  +0.24%  "

# Row 13
$ws.Range("D13").Value = "68.454.92"
$ws.Range("E13").Value = "  -0.01%  "

# Row 14
$ws.Range("E14").Value = "  -2.76%  "

# Row 15
$ws.Range("D15").Value = "'23.27"
$ws.Range("D15").Style = "Normal"  # keep as text, matching original (no numeric coercion)
$ws.Range("E15").Value = "  -1.47%  "

# Row 16
$ws.Range("D16").Value = "'10.56"
$ws.Range("D16").Style = "Normal"  # keep as text, matching original (no numeric coercion)
$ws.Range("E16").Value = "  -4.02%  "

# Row 17
$ws.Range("D17").Value = "'333.37"
$ws.Range("D17").Style = "Normal"  # keep as text, matching original (no numeric coercion)
$ws.Range("E17").Value = "  -2.92%  "

# Row 18
$ws.Range("D18").Value = "'6.87"
$ws.Range("D18").Style = "Normal"  # keep as text, matching original (no numeric coercion)
$ws.Range("E18").Value = "  -4.03%  "

# Row 19
$ws.Range("E19").Value = "  -1.69%  "

# Row 20
$ws.Range("E20").Value = "  -0.01%  "

# Row 21
$ws.Range("E21").Value = "  -0.68%  "

# Row 22
$ws.Range("D22").Value = "'66.25"
$ws.Range("D22").Style = "Normal"  # keep as text, matching original (no numeric coercion)
$ws.Range("E22").Value = "  -2.69%  "

# Row 23
$ws.Range("E23").Value = "  -3.75%  "

# Row 24
$ws.Range("D24").Value = "'8.08"
$ws.Range("D24").Style = "Normal"  # keep as text, matching original (no numeric coercion)
$ws.Range("E24").Value = "  -2.02%  "

# Row 25
$ws.Range("D25").Value = "0.0₃0808"
$ws.Range("E25").Value = "  -3.90%  "

# Row 26
$ws.Range("D26").Value = "'7.14"
$ws.Range("D26").Style = "Normal"  # keep as text, matching original (no numeric coercion)
$ws.Range("E26").Value = "  -2.38%  "

# Row 27
$ws.Range("D27").Value = "'0.999"
$ws.Range("D27").Style = "Normal"  # keep as text, matching original (no numeric coercion)
$ws.Range("E27").Value = "  -0.05%  "

# Row 28
$ws.Range("D28").Value = "'424.27"
$ws.Range("D28").Style = "Normal"  # keep as text, matching original (no numeric coercion)
$ws.Range("E28").Value = "  -2.74%  "

# Row 29
$ws.Range("E29").Value = "  -4.71%  "

# Row 30
$ws.Range("E30").Value = "  -4.88%  "

# Row 31
$ws.Range("E31").Value = "  +0.35%  "

# Row 33
$ws.Range("E33").Value = "  -0.01%  "

# Row 34
$ws.Range("E34").Value = "  -1.24%  "

# Row 35
$ws.Range("D35").Value = "'17.68"
$ws.Range("D35").Style = "Normal"  # keep as text, matching original (no numeric coercion)
$ws.Range("E35").Value = "  -1.36%  "

# Row 36
$ws.Range("E36").Value = "  -2.92%  "

# Row 37
$ws.Range("E37").Value = "  -3.07%  "

# Row 38
$ws.Range("E38").Value = "  -5.47%  "

# Row 39
$ws.Range("E39").Value = "  -2.86%  "

# Row 40
$ws.Range("D40").Value = "'2.03"
$ws.Range("D40").Style = "Normal"  # keep as text, matching original (no numeric coercion)
$ws.Range("E40").Value = "  -2.84%  "

# Row 41
$ws.Range("E41").Value = "  -1.85%  "

# Row 42
$ws.Range("D42").Value = "'128.47"
$ws.Range("D42").Style = "Normal"  # keep as text, matching original (no numeric coercion)
$ws.Range("E42").Value = "  -4.81%  "

# Row 43
$ws.Range("E43").Value = "  -0.43%  "

# Row 44
$ws.Range("E44").Value = "  -1.90%  "

# Row 45
$ws.Range("D45").Value = "'0.558"
$ws.Range("D45").Style = "Normal"  # keep as text, matching original (no numeric coercion)
$ws.Range("E45").Value = "  -0.63%  "

# Row 46
$ws.Range("D46").Value = "'0.0907"
$ws.Range("D46").Style = "Normal"  # keep as text, matching original (no numeric coercion)
$ws.Range("E46").Value = "  -0.80%  "

# Row 47
$ws.Range("E47").Value = "  +0.23%  "

# Row 48
$ws.Range("E48").Value = "  -4.69%  "

# Row 49
$ws.Range("E49").Value = "  -9.58%  "

# Row 50
$ws.Range("D50").Value = "'16.65"
$ws.Range("D50").Style = "Normal"  # keep as text, matching original (no numeric coercion)
$ws.Range("E50").Value = "  -5.94%  "

# Row 51
$ws.Range("B51").Value = "Fantom"
$ws.Range("C51").Value = "https://coinranking.com/coin/uIEWfMFnQo9K_+fantom-ftm"
$ws.Range("D51").Value = "'0.580"
$ws.Range("D51").Style = "Normal"  # keep as text, matching original (no numeric coercion)
$ws.Range("E51").Value = "  -4.57%  "
